$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New s_vals data (B:E, G) for rows 2-6, regenerated to filter save games.
$data = @{
    2 = @{ B = 0.01514828764759746;  C = 0.3127903958511391;  D = 0.1575252929769615;  E = 0.496779210170732;  G = 0.9822431866464301 }
    3 = @{ B = 3.230985683306322;    C = 1.667794583268128;   D = 3.900430680208489;    E = 0.496779210170732;  G = 9.295990156953671 }
    4 = @{ B = 1.459612070389937;    C = 1.667794583268128;   D = 3.900430680208489;    E = 0.496779210170732;  G = 7.524616544037286 }
    5 = @{ B = 0.3048080303191223;   C = 1.667794583268128;   D = 3.900430680208489;    E = 0.496779210170732;  G = 6.369812503966472 }
    6 = @{ B = 3.230985683306322;    C = 1.667794583268128;   D = 0.8054896365839992;   E = 0.496779210170732;  G = 6.201049113329182 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
